$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game")

# Row 17 previously held: A17="doCarTurn()", B17=1, C17="(0,0)", D17=1
# New layout inserts a header row (17) for the new equivalence-class columns,
# pushing the old test-case values down to row 18, and adds a new test case
# (invalid direction -> IllegalArgumentException) on row 19.

# Clear the old B17 / D17 numeric values (their content moves to row 18).
$ws.Range("B17").ClearContents()
$ws.Range("D17").ClearContents()

# Row 17: keep A17 as-is, add the new equivalence-class header in C17.
$ws.Range("C17").Value = "(x,y)"

# Row 18: the former row-17 test case data, plus the new expected-result /
# equivalence-class cells for the acceleration argument check.
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "(0,0)"
$ws.Range("D18").Value = 1
$ws.Range("G18").Value = "two cars; active car is crashed;one car remaining"
$ws.Range("E18").Value = "second car is the winner"

# Row 19: the new invalid-direction test case.
$ws.Range("G19").Value = "invalid direction"
$ws.Range("E19").Value = "IllegalArgumentException"

# New column G needs the same bestFit/custom width treatment as the rest
# (14.1 chars is what this host's pixel-rounding turns into a stored "15").
$ws.Range("G1").ColumnWidth = 14.1

# Scroll/selection state recorded by Excel when the author saved the file.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("E19").Select()
